$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$descPart1 = "yeahLorem ipsum dolor sit amet, consectetur adipiscing elit. Aliquam tellus sem, sodales eget nibh venenatis, dictum laoreet leo. Curabitur eleifend tellus eget tortor consectetur egestas. Nullam finibus pellentesque elit non lacinia. Integer venenatis pellentesque turpis, sit amet faucibus lectus. Duis vitae lacus nisi. Ut maximus nisi quis posuere ultricies. Quisque eu massa ligula. Orci varius natoque penatibus et magnis dis parturient montes, nascetur ridiculus mus. Cras faucibus laoreet nisl non pellentesque. Curabitur volutpat nisi ac varius rutrum. Suspendisse potenti."
$descPart2 = "Aenean efficitur efficitur nisl faucibus scelerisque. Etiam mattis est eu volutpat tristique. Nullam vitae massa consectetur, fermentum elit eu, fringilla nunc. Vivamus ac consectetur purus. Nunc at sem sed odio volutpat tempus. Donec in semper est. Phasellus quis sollicitudin massa. Nullam turpis lorem, eleifend eu laoreet vitae, egestas et ipsum. Curabitur sagittis pharetra blandit. Phasellus posuere augue vitae enim dapibus vulputate."

$descSingleLine = $descPart1 + " " + $descPart2
$descMultiLine = $descPart1 + "`n" + $descPart2

# --- Header row ---
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "title"
$ws.Range("C1").Value = "authors_firstname"
$ws.Range("D1").Value = "authors_surname"
$ws.Range("E1").Value = "description"

# --- Row 2 / Row 3 text values (non-description cells first) ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Pan Tadeusz"
$ws.Range("C2").Value = "Adam"
$ws.Range("D2").Value = "Mickiewicz"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Lalka"
$ws.Range("C3").Value = "Bolesław"
$ws.Range("D3").Value = "Prus"

# The multi-line description is first used on row 3, and the single-line
# variant first appears on row 2 - assign E3 before E2 so the shared
# string table records them in that same order.
$ws.Range("E3").Value = $descMultiLine
$ws.Range("E2").Value = $descSingleLine

# --- Row 4: Que Vadis ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Que Vadis"
$ws.Range("C4").Value = "Henryk"
$ws.Range("D4").Value = "Sienkiewicz "
$ws.Range("E4").Value = $descMultiLine

# --- Row 5: Przedwiośnie ---
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Przedwiośnie"
$ws.Range("C5").Value = "Stefan"
$ws.Range("D5").Value = "Żeromski"
$ws.Range("E5").Value = $descMultiLine

# --- Row 6: Inny Świat ---
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Inny Świat"
$ws.Range("C6").Value = "Gustaw"
$ws.Range("D6").Value = "Herling-Grudziński"
$ws.Range("E6").Value = $descMultiLine

# --- Wrap text + column width for the new description column ---
$ws.Range("E2:E6").WrapText = $true
$ws.Columns("F").ColumnWidth = 15.3

# --- Row heights (left over from manual resize / autofit in the original edit) ---
$ws.Rows("2").RowHeight = 13.5
$ws.Rows("3").RowHeight = 14.5
$ws.Rows("4").RowHeight = 11
$ws.Rows("5").RowHeight = 13
$ws.Rows("6").RowHeight = 10.5

# --- Selection / view ---
$excel.ActiveWindow.Zoom = 100
$ws.Range("E2").Select()
